$d = $word.ActiveDocument

# 1. Patient name
$d.Content.Find.Execute("FABRICIO SANCHEZ BERGAMIN   ", $true, $false, $false, $false, $false, $true, 1, $false, "ADILSON CASSALHO   ", 2) | Out-Null

# 2. Birth date
$d.Content.Find.Execute("10/10/1981   ", $true, $false, $false, $false, $false, $true, 1, $false, "18/09/1980   ", 2) | Out-Null

# 3. Record number
$d.Content.Find.Execute("576", $true, $false, $false, $false, $false, $true, 1, $false, "58825", 2) | Out-Null

# 4. Mother's name
$d.Content.Find.Execute("JUSSARA SANCHEZ BERGAMIN   ", $true, $false, $false, $false, $false, $true, 1, $false, "APARECIDA MADALENA DE OLIVEIRA CASSALHO   ", 2) | Out-Null

# 5. Report date
$d.Content.Find.Execute("22/10/2019   ", $true, $false, $false, $false, $false, $true, 1, $false, "21/10/2017   ", 2) | Out-Null

# 6. Body content: replace RTF blob with plain-text CT report
$old6 = @'
{\rtf1\ansi\ansicpg1252\deff0\deflang1046{\fonttbl{\f0\fnil\fcharset0 }{\f1\fswiss\fcharset0  }{\f2\fswiss\fprq2\fcharset0  }}
{\colortbl ;\red0\green0\blue0;}
\viewkind4\uc1\pard\f0\fs16 .\par
\par
\par
ECOCARDIOGRAMA\par
\par
\par
\cf1\b\f1\fs20  \par
\pard\qj\cf0\b0\f2 R\'cdTMO:\par
Paciente em ritmo card\'edaco regular.\par
\par
C\'c2MARAS CARD\'cdACAS:\par
C\'e2maras card\'edacas com dimens\'f5es normais.\par
\par
VENTR\'cdCULOS:\par
Ventr\'edculo esquerdo apresenta espessura e fun\'e7\'e3o sist\'f3lica preservadas, n\'e3o sendo observadas altera\'e7\'f5es da contra\'e7\'e3o segmentar de parede.\par
An\'e1lise da fun\'e7\'e3o diast\'f3lica do ventr\'edculo esquerdo com padr\'e3o normal.\par
Ventr\'edculo direito apresenta fun\'e7\'e3o sistolica dento da normalidade,\par
\par
V\'c1LVULA MITRAL:\par
Apresenta aspecto e movimenta\'e7\'e3o normais de suas cuspides.\par
O estudo com Doppler e mapeamento de fluxo em cores s\'e3o normais.\par
\par
V\'c1LVULA A\'d3RTICA:\par
Apresenta aspecto e movimenta\'e7\'e3o normais de suas valvulas.\par
O estudo com Doppler e mapeamento de fluxo em cores s\'e3o normais.\par
\par
V\'c1LVULA TRIC\'daSPIDE;\par
Apresenta aspecto e movimenta\'e7\'e3o normais de suas c\'faspides.\par
O estudo com Doppler e mapeamento de fluxo em cores s\'e3o normais.\par
\par
VALVA PULMONAR:\par
Apresenta aspecto e movimenta\'e7\'e3o normais de suas v\'e1lvulas.\par
O estudo com Doppler e mapeamento de fluxo em cores s\'e3o normais.\par
\par
PERIC\'c1RDIO:\par
Peric\'e1rdio com aspecto ecocardiografico normal.\par
\par
AORTA:\par
Seios a\'f3rticos, aorta ascendente e arco a\'f3rtico com dimens\'f5es e fluxos normais.\par
\par
IMPRESS\'c3O DIAGN\'d3STICA:\par
_________________________ \par
\par
\pard ECODOPPLERCARDIOGRAMA DENTRO DA NORMALIDADE PARA O BIOTIPO E FAIXA ET\'c1RIA.\f0\fs16\par
}
_x0000_
'@
$new6 = @'
.
Nome:	Adilson Cassalho			DN: 18/09/1980
Data:	22/10/2017				Convênio: SAMU/SUS
Solicitante: Dra Alexandra O. Somodi
TOMOGRAFIA COMPUTADORIZADA DE ABDOME E PELVE
RELATÓRIO
TÉCNICA
Exame  realizado  com cortes tomográficos computadorizados axiais,  sem a infusão endovenosa de contraste iodado, segundo solicitação do médico assistente.
Salientamos que a não utilização do meio de contraste iodado por via endovenosa prejudica a adequada caracterização das estruturas abdominais.
ANÁLISE
Fígado de topografia, morfologia, situação e dimensões, preservadas, exibindo coeficientes de atenuação homogêneos.
Não há evidência de dilatação das vias biliares intra ou extra-hepáticas, bem como da vesícula biliar.
Baço, pâncreas e adrenais com topografia, dimensões, contornos e densidade normais.
Rins de topografia, morfologia e dimensões preservadas, com coeficientes de atenuação homogêneos, sem a caracterização de cálculos calicinais radiopacos ou hidronefrose.
Aorta e veia cava inferior com posição e calibre normais.
Ausência de linfonodomegalias, líquido livre ou de coleções organizadas na cavidade abdominal.
Bexiga urinária em pequena repleção, com paredes lisas e regulares e conteúdo homogêneo.
Próstata e vesículas seminais sem alterações detectáveis ao método.
OPINIÃO
Tomografia computadorizada do abdome superior e da pelve evidenciando:
Exame sem alterações significativas.
           Dra. Amanda Prist
             CRM-MG: 56.487

'@
$d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $new6, 2) | Out-Null

Write-Output "done"
